$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 5-7 with the surviving data (shifted up from old rows 7, 9, 10)
$ws.Range("A5").Value = "Microscope"
$ws.Range("B5").Value = "INSTRUMENT"
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 0.08333333333333333

$ws.Range("A6").Value = "inductively coupled plasma optical emisssion spectrometer"
$ws.Range("B6").Value = "INSTRUMENTATION"
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 0.08333333333333333

$ws.Range("A7").Value = "stomata"
$ws.Range("B7").Value = "ORGANISM part"
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 0.08333333333333333

# Delete rows 8-10 entirely
$ws.Range("A8:D10").Delete()
